$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.254.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.60%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.520.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.75%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'536.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.13%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'139.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.88%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.22%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -1.46%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.524.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.18%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.67%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.09%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'5.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.47%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.64%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.967.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.46%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'59.207.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.65%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'22.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.73%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +1.36%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.506.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.86%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.39%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.44%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'322.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.13%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.18%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.95%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'62.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.88%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -2.78%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +1.33%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.00%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.25%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.38%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.05%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0₃0766"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.72%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'162.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.78%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.27%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +1.29%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.92%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.75%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.32%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -2.42%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'37.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.38%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.30%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'284.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.15%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'5.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.01%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.32%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.82%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -1.37%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0931"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.14%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'122.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.20%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'18.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.58%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0511"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.63%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.79%  "
$ws.Range("E51").Style = "Normal"
